# Add a new "EQUITY" column (H) to the portfolio table that computes the
# residual allocation not covered by columns E:G (100% - SUM(E:G)), and
# move the active selection to H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("H1").Value = "EQUITY"

# Match the percentage formatting used elsewhere in the sheet (0.00%)
$ws.Range("H2:H5").NumberFormat = "0.00%"

# Residual-equity formulas (these consolidate into an Excel "shared formula"
# group for H3:H5, matching the existing style of the workbook)
$ws.Range("H2").Formula = "=100%-SUM(E2:G2)"
$ws.Range("H3:H5").Formula = "=100%-SUM(E3:G3)"

# Move/restore the selection as recorded in the sheet view
$ws.Range("H10").Select()
